# daily auto push: 2025-10-03 09:26 UTC
# Append the new daily log entry as row 55:
#   A55 = "2025/10/03"  (date text, same literal text as rows 53-54)
#   B55 = "金"          (weekday text, same literal text as rows 53-54)
#   C55 = 17             (hour, number)
#   D55 = 4              (ranking, number)
#
# Columns A/B of the last existing row already contain the exact text we
# need for the new row (2025/10/03 / 金), so rather than typing the date
# string again (which Excel's COM automation would auto-convert into a
# date serial number + new number-format style), we copy the values from
# the previous row. This preserves the text representation exactly and
# does not touch NumberFormat/styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prevRow = $ws.UsedRange.Rows.Count
$newRow = $prevRow + 1

$ws.Cells.Item($prevRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4163)

$ws.Cells.Item($prevRow, 2).Copy()
$ws.Cells.Item($newRow, 2).PasteSpecial(-4163)

$ws.Cells.Item($newRow, 3).Value = 17
$ws.Cells.Item($newRow, 4).Value = 4
